$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 369.9
$ws.Range("I4").Value = 327.66666
$ws.Range("K4").Value = 327.66666
$ws.Range("M4").Value = -213.66666
$ws.Range("H62").Value = 3477.1765
$ws.Range("I62").Value = 2908.8
$ws.Range("K62").Value = 2908.8
$ws.Range("M62").Value = -2284.8
$ws.Range("H65").Value = 3477.1765
$ws.Range("I65").Value = 2908.8
$ws.Range("K65").Value = 14544
$ws.Range("M65").Value = -11424
$ws.Range("H113").Value = 6051.231
$ws.Range("J113").Value = 7993.2
$ws.Range("L113").Value = 7993.2
$ws.Range("N113").Value = -14501.2
$ws.Range("H138").Value = 3577.7778
$ws.Range("I138").Value = 1732.0513
$ws.Range("J138").Value = 6577.0835
$ws.Range("K138").Value = 5196.1539
$ws.Range("L138").Value = 19731.2505
$ws.Range("M138").Value = -56.15390000000025
$ws.Range("N138").Value = -30011.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11907496
$ws.Range("I2").Value = 62501576
$ws.Range("J2").Value = 3006.647
$ws.Range("K2").Value = 62501576
$ws.Range("L2").Value = 3006.647
$ws.Range("M2").Value = -62501463
$ws.Range("N2").Value = -3232.647
$ws.Range("H32").Value = 4349.709
$ws.Range("I32").Value = 2959.884
$ws.Range("K32").Value = 2959.884
$ws.Range("M32").Value = -2672.884
$ws.Range("H61").Value = 2780.9512
$ws.Range("I61").Value = 1249.9375
$ws.Range("J61").Value = 3760.8
$ws.Range("K61").Value = 1249.9375
$ws.Range("L61").Value = 3760.8
$ws.Range("M61").Value = -1037.9375
$ws.Range("N61").Value = -4184.8
$ws.Range("H116").Value = 11907496
$ws.Range("I116").Value = 62501576
$ws.Range("J116").Value = 3006.647
$ws.Range("K116").Value = 62501576
$ws.Range("L116").Value = 3006.647
$ws.Range("M116").Value = -62499282
$ws.Range("N116").Value = -7594.647
$ws.Range("H136").Value = 2780.9512
$ws.Range("I136").Value = 1249.9375
$ws.Range("J136").Value = 3760.8
$ws.Range("K136").Value = 3749.8125
$ws.Range("L136").Value = 11282.4
$ws.Range("M136").Value = -1199.8125
$ws.Range("N136").Value = -16382.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11907496
$ws.Range("I3").Value = 62501576
$ws.Range("J3").Value = 3006.647
$ws.Range("K3").Value = 62501576
$ws.Range("L3").Value = 3006.647
$ws.Range("M3").Value = -62501462
$ws.Range("N3").Value = -3234.647
$ws.Range("H20").Value = 1669.7894
$ws.Range("I20").Value = 831.6
$ws.Range("J20").Value = 2601.111
$ws.Range("K20").Value = 831.6
$ws.Range("L20").Value = 2601.111
$ws.Range("M20").Value = -584.6
$ws.Range("N20").Value = -3095.111
$ws.Range("H94").Value = 544.29034
$ws.Range("I94").Value = 442.66666
$ws.Range("J94").Value = 892.7143
$ws.Range("K94").Value = 442.66666
$ws.Range("L94").Value = 892.7143
$ws.Range("M94").Value = 8.333340000000021
$ws.Range("N94").Value = -1794.7143
$ws.Range("H107").Value = 3124.8
$ws.Range("I107").Value = 1355.5
$ws.Range("J107").Value = 4304.3335
$ws.Range("K107").Value = 1355.5
$ws.Range("L107").Value = 4304.3335
$ws.Range("M107").Value = 564.5
$ws.Range("N107").Value = -8144.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1593.2
$ws.Range("I107").Value = 419.2143
$ws.Range("J107").Value = 2620.4375
$ws.Range("K107").Value = 419.2143
$ws.Range("L107").Value = 2620.4375
$ws.Range("M107").Value = 1500.7857
$ws.Range("N107").Value = -6460.4375
$ws.Range("H122").Value = 3288.25
$ws.Range("I122").Value = 3089.25
$ws.Range("J122").Value = 3487.25
$ws.Range("K122").Value = 9267.75
$ws.Range("L122").Value = 10461.75
$ws.Range("M122").Value = -6817.75
$ws.Range("N122").Value = -15361.75
$ws.Range("H132").Value = 3944.2593
$ws.Range("I132").Value = 2506.9285
$ws.Range("J132").Value = 5492.154
$ws.Range("K132").Value = 7520.7855
$ws.Range("L132").Value = 16476.462
$ws.Range("M132").Value = -4990.7855
$ws.Range("N132").Value = -21536.462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 3777.7778
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -14122
$ws.Range("H77").Value = 3777.7778
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 36000
$ws.Range("N77").Value = -46608
$ws.Range("H131").Value = 1568.7
$ws.Range("I131").Value = 2671.6667
$ws.Range("J131").Value = 1096
$ws.Range("K131").Value = 8015.000100000001
$ws.Range("L131").Value = 3288
$ws.Range("M131").Value = -2975.000100000001
$ws.Range("N131").Value = -13368

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 71433310
$ws.Range("I16").Value = 125002050
$ws.Range("J16").Value = 8334
$ws.Range("K16").Value = 125002050
$ws.Range("L16").Value = 8334
$ws.Range("M16").Value = -125001880
$ws.Range("N16").Value = -8674
$ws.Range("H40").Value = 3750
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3750
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3750
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4022
$ws.Range("H136").Value = 3127987.5
$ws.Range("I136").Value = 5265968.5
$ws.Range("J136").Value = 3246.1538
$ws.Range("K136").Value = 15797905.5
$ws.Range("L136").Value = 9738.4614
$ws.Range("M136").Value = -15795355.5
$ws.Range("N136").Value = -14838.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 986.2
$ws.Range("I81").Value = 627.5
$ws.Range("J81").Value = 1225.3334
$ws.Range("K81").Value = 1255
$ws.Range("L81").Value = 2450.6668
$ws.Range("M81").Value = -194
$ws.Range("N81").Value = -4572.6668
$ws.Range("H84").Value = 986.2
$ws.Range("I84").Value = 627.5
$ws.Range("J84").Value = 1225.3334
$ws.Range("K84").Value = 6275
$ws.Range("L84").Value = 12253.334
$ws.Range("M84").Value = -971
$ws.Range("N84").Value = -22861.334
$ws.Range("H100").Value = 475.73334
$ws.Range("I100").Value = 419.66666
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 839.33332
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -298.33332
$ws.Range("N100").Value = -2482
$ws.Range("H122").Value = 296099.6
$ws.Range("I122").Value = 417995.25
$ws.Range("K122").Value = 1253985.75
$ws.Range("M122").Value = -1253985.75
